# Auto-generated edit script: updates cryptocurrency price/volume/hour data
# per the GitHub Actions symbol-list refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# All D (Price), E (Volume 1h) and G (Hora) columns hold plain text in this
# sheet (t="inlineStr"), even though the strings look numeric. Force the
# "Text" number format before writing so Excel does not silently reinterpret
# the new value as a Number/Percentage (which would change the cell type and
# could drop meaningful trailing zeros, e.g. "0.0001250").
function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "291.55"
Set-TextValue $ws.Range("E2") "-3.00%"
Set-TextValue $ws.Range("G2") "10"

# Row 3
Set-TextValue $ws.Range("D3") "30.78"
Set-TextValue $ws.Range("E3") "-5.46%"
Set-TextValue $ws.Range("G3") "10"

# Row 4
Set-TextValue $ws.Range("E4") "-0.26%"
Set-TextValue $ws.Range("G4") "10"

# Row 5
Set-TextValue $ws.Range("D5") "0.07216"
Set-TextValue $ws.Range("E5") "-5.85%"
Set-TextValue $ws.Range("G5") "10"

# Row 6
Set-TextValue $ws.Range("D6") "1.829"
Set-TextValue $ws.Range("E6") "-5.35%"
Set-TextValue $ws.Range("G6") "10"

# Row 7
Set-TextValue $ws.Range("D7") "7.697"
Set-TextValue $ws.Range("E7") "-1.73%"
Set-TextValue $ws.Range("G7") "10"

# Row 8
Set-TextValue $ws.Range("D8") "3.766"
Set-TextValue $ws.Range("E8") "-0.94%"
Set-TextValue $ws.Range("G8") "10"

# Row 9
Set-TextValue $ws.Range("D9") "0.8977"
Set-TextValue $ws.Range("E9") "-2.25%"
Set-TextValue $ws.Range("G9") "10"

# Row 10
Set-TextValue $ws.Range("D10") "0.1653"
Set-TextValue $ws.Range("E10") "-5.33%"
Set-TextValue $ws.Range("G10") "10"

# Row 11
Set-TextValue $ws.Range("D11") "0.07681"
Set-TextValue $ws.Range("E11") "-1.20%"
Set-TextValue $ws.Range("G11") "10"

# Row 12
Set-TextValue $ws.Range("D12") "0.07960"
Set-TextValue $ws.Range("E12") "-6.66%"
Set-TextValue $ws.Range("G12") "10"

# Row 13
Set-TextValue $ws.Range("D13") "0.03035"
Set-TextValue $ws.Range("E13") "-4.11%"
Set-TextValue $ws.Range("G13") "10"

# Row 14
Set-TextValue $ws.Range("E14") "0.19%"
Set-TextValue $ws.Range("G14") "10"

# Row 15
Set-TextValue $ws.Range("D15") "0.001496"
Set-TextValue $ws.Range("E15") "-1.02%"
Set-TextValue $ws.Range("G15") "10"

# Row 16
Set-TextValue $ws.Range("D16") "0.005690"
Set-TextValue $ws.Range("E16") "-4.19%"
Set-TextValue $ws.Range("G16") "10"

# Row 17
Set-TextValue $ws.Range("G17") "10"

# Row 18
Set-TextValue $ws.Range("D18") "3.469"
Set-TextValue $ws.Range("E18") "0.11%"
Set-TextValue $ws.Range("G18") "10"

# Row 19
Set-TextValue $ws.Range("D19") "2.080"
Set-TextValue $ws.Range("E19") "-3.40%"
Set-TextValue $ws.Range("G19") "10"

# Row 20
Set-TextValue $ws.Range("D20") "0.3319"
Set-TextValue $ws.Range("E20") "-1.05%"
Set-TextValue $ws.Range("G20") "10"

# Row 21
Set-TextValue $ws.Range("D21") "0.1330"
Set-TextValue $ws.Range("E21") "0.29%"
Set-TextValue $ws.Range("G21") "10"

# Row 22
Set-TextValue $ws.Range("D22") "4.049"
Set-TextValue $ws.Range("E22") "-5.39%"
Set-TextValue $ws.Range("G22") "10"

# Row 23
Set-TextValue $ws.Range("D23") "0.2387"
Set-TextValue $ws.Range("E23") "19.83%"
Set-TextValue $ws.Range("G23") "10"

# Row 24
Set-TextValue $ws.Range("D24") "0.04508"
Set-TextValue $ws.Range("E24") "0.10%"
Set-TextValue $ws.Range("G24") "10"

# Row 25
Set-TextValue $ws.Range("E25") "-0.49%"
Set-TextValue $ws.Range("G25") "10"

# Row 26
Set-TextValue $ws.Range("D26") "0.004010"
Set-TextValue $ws.Range("E26") "-8.98%"
Set-TextValue $ws.Range("G26") "10"

# Row 27
Set-TextValue $ws.Range("D27") "0.0001250"
Set-TextValue $ws.Range("E27") "-0.15%"
Set-TextValue $ws.Range("G27") "10"

# Row 28
Set-TextValue $ws.Range("G28") "10"

# Row 29
Set-TextValue $ws.Range("G29") "10"

# Row 30
Set-TextValue $ws.Range("G30") "10"

# Row 31
Set-TextValue $ws.Range("G31") "10"

# Row 32
Set-TextValue $ws.Range("G32") "10"

# Row 33
Set-TextValue $ws.Range("G33") "10"

# Row 34
Set-TextValue $ws.Range("G34") "10"

# Row 35
Set-TextValue $ws.Range("G35") "10"

# Row 36
Set-TextValue $ws.Range("G36") "10"

# Row 37
Set-TextValue $ws.Range("G37") "10"

# Row 38
Set-TextValue $ws.Range("G38") "10"

# Row 39
Set-TextValue $ws.Range("D39") "0.01589"
Set-TextValue $ws.Range("E39") "-6.24%"
Set-TextValue $ws.Range("G39") "10"

# Row 40
Set-TextValue $ws.Range("D40") "0.04399"
Set-TextValue $ws.Range("E40") "-6.10%"
Set-TextValue $ws.Range("G40") "10"

# Row 41
Set-TextValue $ws.Range("D41") "0.007306"
Set-TextValue $ws.Range("E41") "-2.53%"
Set-TextValue $ws.Range("G41") "10"

# Row 42
Set-TextValue $ws.Range("D42") "0.01002"
Set-TextValue $ws.Range("G42") "10"

# Row 43
Set-TextValue $ws.Range("D43") "0.1307"
Set-TextValue $ws.Range("E43") "-3.26%"
Set-TextValue $ws.Range("G43") "10"

# Row 44
Set-TextValue $ws.Range("D44") "0.002006"
Set-TextValue $ws.Range("E44") "-14.02%"
Set-TextValue $ws.Range("G44") "10"

# Row 45
Set-TextValue $ws.Range("D45") "0.009505"
Set-TextValue $ws.Range("E45") "-9.72%"
Set-TextValue $ws.Range("G45") "10"

# Row 46
Set-TextValue $ws.Range("D46") "0.00005927"
Set-TextValue $ws.Range("E46") "-5.30%"
Set-TextValue $ws.Range("G46") "10"

# Row 47
Set-TextValue $ws.Range("E47") "-0.10%"
Set-TextValue $ws.Range("G47") "10"

# Row 48
Set-TextValue $ws.Range("G48") "10"

# Row 49
Set-TextValue $ws.Range("G49") "10"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "-0.10%"
Set-TextValue $ws.Range("G50") "10"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.10%"
Set-TextValue $ws.Range("G51") "10"

